$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.736.38"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").Value = "1.646.68"
$ws.Range("E3").Value = "  +0.63%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.505"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("E7").Value = "  +0.42%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.251"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.45%  "

$ws.Range("E9").Value = "  -0.18%  "

$ws.Range("E10").Value = "  +1.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0843"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("D12").Value = "1.878.86"
$ws.Range("E12").Value = "  +0.85%  "

$ws.Range("D13").Value = "1.644.08"
$ws.Range("E13").Value = "  +0.84%  "

$ws.Range("E14").Value = "  +0.95%  "

$ws.Range("E15").Value = "  +1.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.25"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = "26.764.35"
$ws.Range("E17").Value = "  +0.69%  "

$ws.Range("D18").Value = "0.0₃0743"
$ws.Range("E18").Value = "  -0.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.77%  "

$ws.Range("E20").Value = "  +0.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +13.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.95%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.56%  "

$ws.Range("E26").Value = "  +0.48%  "

$ws.Range("E27").Value = "  -0.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.45%  "

$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("E31").Value = "  +1.28%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.45%  "

$ws.Range("E33").Value = "  +0.79%  "

$ws.Range("D34").Value = "1.277.59"
$ws.Range("E34").Value = "  +1.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.17%  "

$ws.Range("E36").Value = "  +2.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0178"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.85%  "

$ws.Range("E38").Value = "  +5.61%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.829"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.69%  "

$ws.Range("E40").Value = "  +0.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.813"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.88%  "

$ws.Range("E42").Value = "  -1.42%  "

$ws.Range("E43").Value = "  +1.46%  "

$ws.Range("D44").Value = "1.790.40"
$ws.Range("E44").Value = "  +1.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "91.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.50%  "

$ws.Range("E47").Value = "  +0.57%  "

$ws.Range("E48").Value = "  +1.08%  "

$ws.Range("E49").Value = "  +0.77%  "

$ws.Range("E50").Value = "  +2.57%  "

$ws.Range("E51").Value = "  +1.73%  "

